# Generate Report for Handoff
# - Remove the "5f2c9c7b-..." row (row 3) from all three sheets.
# - Update the "0e314636-..." row (row 2) status/date values to reflect the
#   new handoff ("Ready for handoff") instead of the old handback status.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Update row 2 values first (row 3 is about to be deleted).
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-22 08:47:39"

# Drop every hyperlink on the sheet, then re-create only the one that must
# survive (the row-3 hyperlink will not be re-added).
$overview.Hyperlinks.Delete()

# Remove the now-obsolete row for 5f2c9c7b-...
$overview.Range("A3").EntireRow.Delete()

$overview.Hyperlinks.Add($overview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fa37aeda7d79b054fa69d8892cd2d5030b19c0aa/e2e/0e314636-8642-429d-95e2-56fccc4a9f14.md", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-22 08:47:36"

$zhcn.Hyperlinks.Delete()

$zhcn.Range("A3").EntireRow.Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fa37aeda7d79b054fa69d8892cd2d5030b19c0aa/e2e/0e314636-8642-429d-95e2-56fccc4a9f14.md", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd32137e0156987723ec48528e74e08fbb6a24b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.zh-cn.xlf", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/082d98b748f010ed2ff4ac7e07a9fff524fc1d4a/e2e/0e314636-8642-429d-95e2-56fccc4a9f14.md", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0e15d557fcd317b4fe41ec1f8db40daf011aa735/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.zh-cn.xlf", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-22 08:47:39"

$dede.Hyperlinks.Delete()

$dede.Range("A3").EntireRow.Delete()

$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fa37aeda7d79b054fa69d8892cd2d5030b19c0aa/e2e/0e314636-8642-429d-95e2-56fccc4a9f14.md", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dc8113d84e786ed2343e4f690c9ec11183cc6b2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.de-de.xlf", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ea4ef1dac429c9cc4292b87ba8f28c0d4089476d/e2e/0e314636-8642-429d-95e2-56fccc4a9f14.md", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ae94280c9dfee700f60eb19b862f19c3784c13a5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.de-de.xlf", "", "", "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.de-de.xlf") | Out-Null
